# Applies the "Updated symbol list" refresh: new Price (column D) and
# Volume(1h) (column E) figures for the affected coin rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextCell($ref, $text) {
    $cell = $ws.Range($ref)
    # Force text storage so numeric-looking strings (prices, percents)
    # are not auto-converted to numbers, then restore the default style
    # so no stray cell formatting is introduced.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextCell "D2" "261.06"
Set-TextCell "E2" "0.02%"
Set-TextCell "D3" "26.94"
Set-TextCell "E3" "-1.14%"
Set-TextCell "D4" "4.731"
Set-TextCell "E4" "-0.15%"
Set-TextCell "D5" "0.06220"
Set-TextCell "E5" "2.23%"
Set-TextCell "D6" "6.743"
Set-TextCell "E6" "1.13%"
Set-TextCell "D7" "0.8499"
Set-TextCell "E7" "0.57%"
Set-TextCell "D8" "0.9129"
Set-TextCell "E8" "-0.99%"
Set-TextCell "D9" "0.1405"
Set-TextCell "E9" "-0.01%"
Set-TextCell "D10" "0.04998"
Set-TextCell "E10" "0.64%"
Set-TextCell "D11" "0.07076"
Set-TextCell "E11" "-0.31%"
Set-TextCell "D12" "0.03102"
Set-TextCell "E12" "-1.02%"
Set-TextCell "D13" "0.09056"
Set-TextCell "E13" "-0.23%"
Set-TextCell "D14" "0.001527"
Set-TextCell "E14" "-0.37%"
Set-TextCell "D15" "0.0006177"
Set-TextCell "E15" "1.68%"
Set-TextCell "D16" "0.005962"
Set-TextCell "E16" "-3.55%"
Set-TextCell "E17" "-0.27%"
Set-TextCell "D18" "3.171"
Set-TextCell "E18" "0.74%"
Set-TextCell "E19" "-1.04%"
Set-TextCell "E20" "-0.68%"
Set-TextCell "E21" "1.09%"
Set-TextCell "D22" "4.100"
Set-TextCell "E22" "0.16%"
Set-TextCell "D23" "0.04226"
Set-TextCell "E23" "-0.52%"
Set-TextCell "D24" "0.001203"
Set-TextCell "E24" "-1.47%"
Set-TextCell "D25" "0.004076"
Set-TextCell "E25" "4.20%"
Set-TextCell "E26" "0.09%"
Set-TextCell "D40" "0.03949"
Set-TextCell "E40" "2.00%"
Set-TextCell "D41" "0.1112"
Set-TextCell "E41" "-0.09%"
Set-TextCell "D42" "0.004140"
Set-TextCell "E42" "0.24%"
Set-TextCell "D43" "0.002212"
Set-TextCell "E43" "0.19%"
Set-TextCell "D44" "0.01338"
Set-TextCell "E44" "-18.14%"
Set-TextCell "D45" "0.00005164"
Set-TextCell "E45" "-2.89%"
Set-TextCell "D46" "0.00000000751"
Set-TextCell "E46" "0.09%"
Set-TextCell "E47" "-37.55%"
Set-TextCell "D48" "0.2519"
Set-TextCell "E48" "86.15%"
Set-TextCell "D49" "0.00002102"
Set-TextCell "E49" "0.09%"
Set-TextCell "D50" "0.0002002"
Set-TextCell "E50" "0.09%"
